# Add 2022-Q4 data
#
# The workbook currently has sheets: 总计, 2022-Q2, 2022-Q1, 2021-Q4, 2021-Q3.
# A new quarter (2022-Q4) of data has arrived. The "2022-Q2" sheet slot is
# reused to hold the freshest quarter's numbers (2022-Q4) and is renamed;
# a duplicate sheet is inserted right after it so the original "2022-Q2"
# figures are preserved under their own tab. The summary sheet ("总计")
# gets a new leading row for 2022-Q4 and every other row shifts down one.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Worksheets: duplicate "2022-Q2" so its original data survives under
#    its own tab, then repurpose the original tab for the new quarter.
# ---------------------------------------------------------------------
$q2 = $wb.Worksheets.Item("2022-Q2")
$q2.Copy($null, $q2)

# Original sheet (still in slot 2, still holding the old 2022-Q2 numbers)
# becomes the new "2022-Q4" tab.
$q4 = $wb.Worksheets.Item(2)
$q4.Name = "2022-Q4"

# The freshly made copy (slot 3) keeps the old figures, so it becomes the
# "2022-Q2" tab again.
$q2copy = $wb.Worksheets.Item(3)
$q2copy.Name = "2022-Q2"

# Update the fund figures on the "2022-Q4" tab to the new quarter's values.
# Columns D-G are stored as text (matching the source data format); H is numeric.
$q4.Range("D2").Value = "'12.63"
$q4.Range("D2").Style = "Normal"
$q4.Range("E2").Value = "'60.37"
$q4.Range("E2").Style = "Normal"
$q4.Range("F2").Value = "'1.93"
$q4.Range("F2").Style = "Normal"
$q4.Range("G2").Value = "'0.2438"
$q4.Range("G2").Style = "Normal"
$q4.Range("H2").Value = 10

# ---------------------------------------------------------------------
# 2) Summary sheet ("总计"): insert a new leading data row for 2022-Q4 and
#    push the existing quarters down by one row.
# ---------------------------------------------------------------------
$totals = $wb.Worksheets.Item("总计")

# Make room for the new row 6 (shifted-down 2021-Q3 row), copying the
# index-column formatting that the other rows already use.
$totals.Range("A5").Copy()
$totals.Range("A6").PasteSpecial(-4122)
$totals.Range("A6").Value = 4

# Shift rows 2-5 down to rows 3-6, bottom-up so we never overwrite a
# source cell before it has been read.
$totals.Range("B6").Value = $totals.Range("B5").Value2
$totals.Range("C6").Value = $totals.Range("C5").Value2
$totals.Range("D6").Value = $totals.Range("D5").Value2

$totals.Range("B5").Value = $totals.Range("B4").Value2
$totals.Range("C5").Value = $totals.Range("C4").Value2
$totals.Range("D5").Value = $totals.Range("D4").Value2

$totals.Range("B4").Value = $totals.Range("B3").Value2
$totals.Range("C4").Value = $totals.Range("C3").Value2
$totals.Range("D4").Value = $totals.Range("D3").Value2

$totals.Range("B3").Value = $totals.Range("B2").Value2
$totals.Range("C3").Value = $totals.Range("C2").Value2
$totals.Range("D3").Value = $totals.Range("D2").Value2

# New first data row: 2022-Q4
$totals.Range("B2").Value = "2022-Q4"
$totals.Range("C2").Value = 1
$totals.Range("D2").Value = 0.24
